# Add a "price per meter" (ppm) column in column I.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, styled like the other header cells (B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "ppm"

# Data rows: meters are in column G, price(kinda) is in column H.
# ppm = price / meters, rounded to 4 decimal places.
for ($row = 2; $row -le 25; $row++) {
    $meters = $ws.Cells.Item($row, 7).Value2
    $price = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = [Math]::Round($price / $meters, 4)
}
